$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.652.34"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.016.23"
$ws.Range("E3").Value = "  -4.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.43"
$ws.Range("E5").Value = "  -3.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5054"
$ws.Range("E7").Value = "  -3.40%  "
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.05"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09185"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  -3.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.64"
$ws.Range("E12").Value = "  -6.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.163"
$ws.Range("E13").Value = "  -6.03%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.570"
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.977.98"
$ws.Range("E15").Value = "  -6.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.27"
$ws.Range("E16").Value = "  -5.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001128"
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06656"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.98"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.017"
$ws.Range("D23").Value = "29.688.06"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -4.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.35"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  -5.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.588"
$ws.Range("E28").Value = "  -2.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.356"
$ws.Range("E29").Value = "  -7.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.02"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("E31").Value = "  -7.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.603"
$ws.Range("E32").Value = "  -8.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09993"
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.902"
$ws.Range("E34").Value = "  -5.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.805"
$ws.Range("E35").Value = "  -3.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.699"
$ws.Range("E36").Value = "  -8.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02485"
$ws.Range("E37").Value = "  -5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.313"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06415"
$ws.Range("E39").Value = "  -6.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6594"
$ws.Range("E40").Value = "  -6.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.85"
$ws.Range("E42").Value = "  -6.56%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6374"
$ws.Range("E44").Value = "  -6.90%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.66"
$ws.Range("E45").Value = "  -6.45%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.229"
$ws.Range("E46").Value = "  -5.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.293"
$ws.Range("E47").Value = "  -7.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.537"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.144"
$ws.Range("E50").Value = "  -4.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000322"
$ws.Range("E51").Value = "  -6.35%  "
